$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.911.14"
$ws.Range("E2").Value = "  -4.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.218.46"
$ws.Range("E3").Value = "  -6.42%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.38"
$ws.Range("E5").Value = "  +1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.18"
$ws.Range("E6").Value = "  -8.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("E7").Value = "  -6.58%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  -8.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.75"
$ws.Range("E10").Value = "  -10.07%  "

$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0828"
$ws.Range("E12").Value = "  -9.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.76"
$ws.Range("E13").Value = "  -8.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  -3.43%  "

# Swap rows 15 and 16 (Polygon now ranks above WrappedliquidstakedEther2.0)
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.861"
$ws.Range("E15").Value = "  -11.75%  "

$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.554.00"
$ws.Range("E16").Value = "  -6.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("E17").Value = "  -6.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.210.74"
$ws.Range("E18").Value = "  -7.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.836.35"
$ws.Range("E19").Value = "  -5.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  +4.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -9.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  -11.43%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.36"
$ws.Range("E23").Value = "  -10.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.17"
$ws.Range("E24").Value = "  -9.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "236.49"
$ws.Range("E25").Value = "  -8.86%  "

$ws.Range("E26").Value = "  -8.74%  "

$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("E28").Value = "  -9.74%  "

$ws.Range("E29").Value = "  -5.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("E30").Value = "  -12.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.48"
$ws.Range("E31").Value = "  -8.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0880"
$ws.Range("E32").Value = "  -8.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.15"
$ws.Range("E33").Value = "  -8.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.95"
$ws.Range("E34").Value = "  -7.23%  "

$ws.Range("E35").Value = "  -6.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.20"
$ws.Range("E36").Value = "  +9.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("E37").Value = "  +11.51%  "

$ws.Range("E38").Value = "  -6.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.45"
$ws.Range("E39").Value = "  -6.07%  "

# Swap rows 40 and 41 (NEARProtocol now ranks above Kaspa)
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.75"
$ws.Range("E40").Value = "  -4.29%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.103"
$ws.Range("E41").Value = "  -11.76%  "

$ws.Range("E42").Value = "  -8.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.895.79"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.29"
$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.10"
$ws.Range("E46").Value = "  -10.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.207"
$ws.Range("E47").Value = "  -9.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("E48").Value = "  -3.78%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "60.63"
$ws.Range("E49").Value = "  -13.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.47"
$ws.Range("E50").Value = "  -7.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.861"
$ws.Range("E51").Value = "  +15.51%  "
